# Apply the "clean_values" sheet edit: insert 3 new lookup rows
# (boolean TRUE/"Yes" row, "google"->"Google" row, and 3 cloud-platform rows),
# re-point the view's selection, and set the page to portrait orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("clean_values")
$ws.Activate()

# --- Make room for the new rows -------------------------------------------------
# Old layout:
#   15 (A=blank,B=N,C=No) 16 (B=yes,C=Yes) 17 (B=no,C=No)
#   18 (A=cloud_plaftforms,B=Amazon Web Services,C=AWS)
#   19 (A=blank,B=Google Platform,C=Google) 20 (A=blank,B=Microsoft Azure,C=Azure)
# New layout inserts one blank row before old row 15 (for the new TRUE/Yes row)
# and a second blank row after the (shifted) Google Platform row (for the new
# google/Google row); 3 more rows are appended at the end for ibm/IBM,
# "Digital ocean"/"Digital Ocean" and "oracle, AWS"/"Oracle, AWS".
$ws.Rows.Item(15).Insert()
$ws.Rows.Item(21).Insert()

# --- Row 15 (new): boolean TRUE formatted as left-aligned text, "Yes" ----------
$ws.Range("B15").Value = $true
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").HorizontalAlignment = -4131
$ws.Range("C15").Value = "Yes"

# --- Row 21 (new): google -> Google --------------------------------------------
$ws.Range("B21").Value = "google"
$ws.Range("C21").Value = "Google"

# --- New rows appended at the bottom (23-25), styled with a distinct font -------
$ws.Range("B23").Value = "ibm"
$ws.Range("C23").Value = "IBM"
$ws.Range("B24").Value = "Digital ocean"
$ws.Range("C24").Value = "Digital Ocean"
$ws.Range("B25").Value = "oracle, AWS"
$ws.Range("C25").Value = "Oracle, AWS"
$ws.Range("B23:C25").Font.Name = "Calibri "

# --- View state: scroll so row 7 is at the top, select H18 ---------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("H18").Select()

# --- Page setup: portrait orientation -------------------------------------------
$ws.PageSetup.Orientation = 1
